$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values for specific rows as per repulled data
$ws.Range("F2").Value = 1
$ws.Range("F5").Value = -10
$ws.Range("F7").Value = -5
$ws.Range("F12").Value = -3
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = 2
